$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.30"
$ws.Range("E2").Value = "'1.89%"
$ws.Range("D3").Value = "'27.31"
$ws.Range("E3").Value = "'1.71%"
$ws.Range("D4").Value = "'4.705"
$ws.Range("E4").Value = "'0.56%"
$ws.Range("D5").Value = "'0.06094"
$ws.Range("E5").Value = "'3.34%"
$ws.Range("D6").Value = "'6.673"
$ws.Range("E6").Value = "'1.08%"
$ws.Range("D7").Value = "'0.8459"
$ws.Range("E7").Value = "'-0.60%"
$ws.Range("D8").Value = "'0.9243"
$ws.Range("E8").Value = "'0.48%"
$ws.Range("D9").Value = "'0.1407"
$ws.Range("E9").Value = "'2.09%"
$ws.Range("D10").Value = "'0.04705"
$ws.Range("E10").Value = "'13.09%"
$ws.Range("D11").Value = "'0.07106"
$ws.Range("E11").Value = "'1.31%"
$ws.Range("D12").Value = "'0.03103"
$ws.Range("E12").Value = "'1.68%"
$ws.Range("D13").Value = "'0.09064"
$ws.Range("E13").Value = "'-0.47%"
$ws.Range("D14").Value = "'0.001543"
$ws.Range("E14").Value = "'0.98%"
$ws.Range("D15").Value = "'0.0006062"
$ws.Range("E15").Value = "'0.04%"
$ws.Range("D16").Value = "'0.006146"
$ws.Range("E16").Value = "'0.00%"
$ws.Range("D17").Value = "'3.448"
$ws.Range("E17").Value = "'-0.63%"
$ws.Range("E18").Value = "'-0.76%"
$ws.Range("E21").Value = "'-0.52%"
$ws.Range("D22").Value = "'4.084"
$ws.Range("E22").Value = "'4.68%"
$ws.Range("E23").Value = "'0.05%"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'0.03%"
$ws.Range("E25").Value = "'-8.42%"
$ws.Range("E26").Value = "'0.08%"
$ws.Range("E27").Value = "'3.39%"
$ws.Range("D40").Value = "'0.03876"
$ws.Range("E40").Value = "'2.50%"
$ws.Range("E41").Value = "'1.49%"
$ws.Range("D42").Value = "'0.004103"
$ws.Range("E42").Value = "'6.11%"
$ws.Range("D43").Value = "'0.01628"
$ws.Range("E43").Value = "'15.01%"
$ws.Range("E44").Value = "'-7.21%"
$ws.Range("D45").Value = "'0.00005157"
$ws.Range("E45").Value = "'-3.65%"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("E47").Value = "'19.67%"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("E50").Value = "'0.03%"
